# The header row in each worksheet currently reads (columns C-F):
#   C1=normalize_group, D1=trajgroup_no_vary_q, E1=uniform_scaling_q, F1=variable_trajectory_group
# It needs to become:
#   C1=variable_trajectory_group, D1=normalize_group, E1=trajgroup_no_vary_q, F1=uniform_scaling_q
# i.e. "variable_trajectory_group" moves from F1 to C1, and the other three
# headers each shift one column to the right (C->D, D->E, E->F).
# All data rows below the header have empty cells in columns C-G, so only the
# header row needs to be touched.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Item(1, 3).Value = "variable_trajectory_group"
    $ws.Cells.Item(1, 4).Value = "normalize_group"
    $ws.Cells.Item(1, 5).Value = "trajgroup_no_vary_q"
    $ws.Cells.Item(1, 6).Value = "uniform_scaling_q"
}
